$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column J holds "mb two step" results; add header then values row by row,
# and fill the new D-column values for rows 23-25 (D22 already populated).
# Values are entered row-by-row (D before J) to mirror natural authoring order.

$ws.Range("J21").Value = "mb two step"

$ws.Range("J22").Value = "0.20(0.026)"

$ws.Range("D23").Value = "0.59(0.009)"
$ws.Range("J23").Value = "0.61(0.009)"

$ws.Range("D24").Value = "0.63(0.011)"
$ws.Range("J24").Value = "0.60(0.010)"

$ws.Range("D25").Value = "0.75(0.009)"
$ws.Range("J25").Value = "0.75(0.010)"

# Update the selection to match the new active cell/range noted in the diff
$ws.Range("D26").Select()
